# Zadania.docx edit:
# In the exercise answer paragraph "Odp." / "32" / "36" / "38" / "52"
# (the answers to the "divisible by 2" exercise), the line break and the
# value "32" that follow "Odp." used to live in the same run. After the
# edit they are split into two separate runs: one containing just the
# manual line break, and one containing just the text "32" - matching how
# the other "Odp." answer lists in the document are already structured
# (each number on its own run following its own break).

$d = $word.ActiveDocument

# Locate the paragraph robustly by searching for the literal text
# "Odp.<line break>32" (a vertical-tab character, Chr(11)/^l, is how a
# manual line break - <w:br/> - shows up in Range.Text).
$vt = [char]11
$needle = "Odp." + $vt + "32"

$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.MatchWildcards = $false
$find.Text = $needle
$found = $find.Execute()

if ($found) {
    # $rng now spans "Odp.<br>32" - the trailing 2 characters are "32".
    $numberRange = $d.Range($rng.End - 2, $rng.End)

    # Toggling a character formatting property and then reverting it
    # forces Word to give this text its own run, splitting it off from
    # the run that holds the preceding manual line break - without
    # altering the visible formatting of the text.
    $numberRange.Font.Bold = $true
    $numberRange.Font.Bold = $false

    Write-Host "Split 'Odp./32' run at [$($numberRange.Start)-$($numberRange.End)]: '$($numberRange.Text)'"
} else {
    Write-Host "Target text 'Odp.<br>32' not found - no changes made."
}

$d.Save()
